$d = $word.ActiveDocument

# 1. "Check the dataset " + "size " + "to get an idea..." -> single run
$d.Content.Find.Execute(
    "Check the dataset size to get an idea of the data volume and potential computational requirements.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Check the dataset size to get an idea of the data volume and potential computational requirements.",
    2) | Out-Null

# 2. "Step 2:" (bold, untouched) + " " + "Check for Data Quality:" -> "Step 2:" + " Check for Data Quality:"
$d.Content.Find.Execute(
    " Check for Data Quality:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Check for Data Quality:",
    2) | Out-Null

# 3. "Identify missing values" + " or " + "inconsistencies in the data." -> single run
$d.Content.Find.Execute(
    "Identify missing values or inconsistencies in the data.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Identify missing values or inconsistencies in the data.",
    2) | Out-Null

# 4. "I" + "dentify outliers" (paragraph) then separate paragraph "." -> merge into one
#    paragraph "Identify outliers." (paragraph mark between them is removed).
$d.Content.Find.Execute(
    "Identify outliers^p.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Identify outliers.",
    2) | Out-Null

# 5. "Step 3:" (bold, untouched) + " " + "Define the Project Goal and Questions:" -> merge
$d.Content.Find.Execute(
    " Define the Project Goal and Questions:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Define the Project Goal and Questions:",
    2) | Out-Null

# 6. "S" + "pecific questions you need to answer." -> single run
$d.Content.Find.Execute(
    "Specific questions you need to answer.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Specific questions you need to answer.",
    2) | Out-Null

# 7. "Step 4:" (bold, untouched) + " " + "Data Cleaning and Preprocessing:" -> merge
$d.Content.Find.Execute(
    " Data Cleaning and Preprocessing:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Data Cleaning and Preprocessing:",
    2) | Out-Null

# 8. "Handle missing values, " + (proofErr spellStart) + "duplicates" + (proofErr spellEnd)
#    -> single run "Handle missing values, duplicates" with no spell-check markers.
$d.Content.Find.Execute(
    "Handle missing values, duplicates",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Handle missing values, duplicates",
    2) | Out-Null

# 9. "Step" + " " + "6:" (bold, untouched) + " " + "Interpret Results:" -> merge last two runs
$d.Content.Find.Execute(
    " Interpret Results:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Interpret Results:",
    2) | Out-Null

# 10. "Step 7:" (bold, untouched) + " " + "Validate Results:" -> merge
$d.Content.Find.Execute(
    " Validate Results:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " Validate Results:",
    2) | Out-Null

# 11. "Perform" + "ed" + " sanity" -> single run "Performed sanity"
$d.Content.Find.Execute(
    "Performed sanity",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Performed sanity",
    2) | Out-Null
